# Add a new bullet paragraph right before the closing "And so on..." line on
# the "What will I do" slide (slide 21), describing the error-handling work
# item. The new paragraph is built out of four runs so that the "etc" run
# keeps its own (flagged) run properties, matching how PowerPoint splits a
# typed sentence into separate runs around an autocorrect/spell-check marker.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(21)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Locate the existing "And so on..." paragraph so the new one can be
# inserted immediately before it, regardless of its current index.
$targetIndex = 0
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    if ($tr.Paragraphs($i).Text -like "And so on*") {
        $targetIndex = $i
        break
    }
}

$lastPara = $tr.Paragraphs($targetIndex)

# Create the new paragraph (first run) by inserting text + a paragraph break
# before the "And so on..." paragraph.
[void]$lastPara.InsertBefore("Add error handling(API key expired, csv `r")

# The freshly created paragraph is now at $targetIndex; it currently holds
# just the first run. Append the remaining three runs to it.
$newPara = $tr.Paragraphs($targetIndex)
[void]$newPara.InsertAfter("syntax error, ")
[void]$newPara.InsertAfter("etc")
[void]$newPara.InsertAfter(")")
